$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999758593883,
    0.99771098623859866,
    0.98907929634546798,
    0.99231962491305825,
    0.9824262207717267,
    0.95856686067711938,
    0.95761960174501382,
    0.95556622172429195,
    0.95353252754346207,
    0.9524352960478375,
    0.95232577386547312,
    0.9523738802789079,
    0.95391371497639144,
    0.95576327049479937,
    0.95859101855179407,
    0.96239030103470413,
    0.97313720719840302,
    0.97202829399716228,
    0.99665772022135468,
    0.98954052633357126,
    0.98814200611991343,
    0.9868774951822552,
    0.97691099290309957,
    0.96389000999740959,
    0.9574329864915283,
    0.95491297286641408,
    0.95006852821778032,
    0.92860258037769827,
    0.91333354170752035,
    0.9067635429863552,
    0.89910986226204259,
    0.89743048813093118,
    0.89691045898095123
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
